$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tweak the two instruction sentences
$ws.Range("A2").Value = "We will first practice listening to a sentence without noise. Press 'Return' to continue."
$ws.Range("A3").Value = "Great! Now we will practice listening to sentences with noise. This may be more challenging, so just try your best! Press 'Return' to continue."

# Add two new practice rows (noisy-sentence items), with a new reminder sentence
$ws.Range("A4").Value = "If you are not sure of what you hear, just do your best. Please guess if possible."
$ws.Range("A5").Value = "If you are not sure of what you hear, just do your best. Please guess if possible."

$ws.Range("B5").Value = "soundfiles/NU1109_0825_talker01_SNR-5.wav"
$ws.Range("B4").Value = "soundfiles/NU1109_0792_talker01_SNR-5.wav"

# Widen the columns to fit the longer text / file names
$ws.Columns.Item(1).ColumnWidth = 149.16666666666666
$ws.Columns.Item(2).ColumnWidth = 40.498697916666664

$ws.Range("A5").Select()
